$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to retain its literal text formatting (e.g. trailing
# zeros / thousand-dot separators) instead of being auto-coerced to a number.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.003.62"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.909.51"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "0.7940"
$ws.Range("E5").Value = "  +5.36%  "
$ws.Range("D6").Value = "242.10"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.3165"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "26.37"
$ws.Range("E9").Value = "  +4.49%  "
$ws.Range("D10").Value = "0.06934"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "0.08000"
$ws.Range("D12").Value = "1.909.75"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "0.7437"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "5.196"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "93.02"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "30.003.98"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "5.875"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("D19").Value = "246.41"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").Value = "0.000007758"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "2.151.79"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "6.845"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "168.32"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "9.249"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "0.1398"
$ws.Range("E27").Value = "  +8.02%  "
$ws.Range("D28").Value = "18.92"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "2.034"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "1.370"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "1.513"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "4.318"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "4.093"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "0.05564"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "0.7332"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "2.787"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "6.126"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "0.4423"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "72.44"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "0.8338"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "1.880"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").Value = "100.61"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "7.541"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "987.51"
$ws.Range("E48").Value = "  +8.04%  "
$ws.Range("D49").Value = "2.057.26"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "36.28"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "2.816"
$ws.Range("E51").Value = "  +7.31%  "

# Drop the temporary text number-format so the cells fall back to the workbook
# default style (matches the source file, which carries no explicit style index
# on these cells).
$priceRange.ClearFormats()
